$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.093.38"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
$ws.Range("D3").Value = "3.256.12"
$ws.Range("E3").Value = "  -0.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.67"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.93%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.97%  "

# Row 9
$ws.Range("E9").Value = "  -1.63%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.06%  "

# Row 11
$ws.Range("E11").Value = "  +0.30%  "

# Row 12
$ws.Range("D12").Value = "3.823.77"
$ws.Range("E12").Value = "  -0.68%  "

# Row 13
$ws.Range("E13").Value = "  -0.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.13"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.95%  "

# Row 15
$ws.Range("D15").Value = "68.192.69"
$ws.Range("E15").Value = "  +0.48%  "

# Row 16
$ws.Range("E16").Value = "  +0.62%  "

# Row 17
$ws.Range("D17").Value = "3.290.94"
$ws.Range("E17").Value = "  +0.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.84"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.43%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "393.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.14%  "

# Row 22
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.47"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.18%  "

# Row 23
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("E24").Value = "  +0.74%  "

# Row 25
$ws.Range("E25").Value = "  -0.76%  "

# Row 26
$ws.Range("E26").Value = "  +4.30%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.73"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("E28").Value = "  -0.03%  "

# Row 29
$ws.Range("E29").Value = "  -0.21%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.69"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.90"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.29%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.14"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.40%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.56"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.45%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.50"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.84%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.69%  "

# Row 38
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.821"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.93%  "

# Row 39
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.79"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.85%  "

# Row 40
$ws.Range("E40").Value = "  -1.03%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.52"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.80%  "

# Row 42
$ws.Range("E42").Value = "  -6.66%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0689"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.08%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.43"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.79%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.654.33"
$ws.Range("E45").Value = "  -0.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.10"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "339.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.56%  "

# Row 48
$ws.Range("E48").Value = "  -0.99%  "

# Row 49
$ws.Range("E49").Value = "  +3.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.49"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.20%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.102"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.41%  "
